$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.900673960402175
$ws.Range("D2").Value = 6.106787525169712
$ws.Range("E2").Value = 13.08979260038256
$ws.Range("F2").Value = 32.17575759233443
$ws.Range("G2").Value = 42.77752766999639
$ws.Range("H2").Value = 17.46021584234498
$ws.Range("L2").Value = 9.524541386855931
$ws.Range("M2").Value = 60.06208870055031
$ws.Range("C3").Value = 10.08891691546194
$ws.Range("D3").Value = 6.26695021523622
$ws.Range("E3").Value = 12.83695780203309
$ws.Range("F3").Value = 32.58558489583338
$ws.Range("G3").Value = 43.40942341121801
$ws.Range("H3").Value = 17.72293429343335
$ws.Range("L3").Value = 9.324046764481693
$ws.Range("M3").Value = 56.78537599877266
$ws.Range("C4").Value = 10.21040316096158
$ws.Range("D4").Value = 6.368535370690585
$ws.Range("E4").Value = 12.68542481973064
$ws.Range("F4").Value = 32.8683713238417
$ws.Range("G4").Value = 43.8532623557407
$ws.Range("H4").Value = 17.89530141460441
$ws.Range("L4").Value = 9.201654695173843
$ws.Range("M4").Value = 54.6705423302256
$ws.Range("C5").Value = 10.26137685743486
$ws.Range("D5").Value = 6.410759512991618
$ws.Range("E5").Value = 12.62467581225446
$ws.Range("F5").Value = 32.99117785101921
$ws.Range("G5").Value = 44.04758060477004
$ws.Range("H5").Value = 17.96826817603767
$ws.Range("L5").Value = 9.152017650203224
$ws.Range("M5").Value = 53.78326798305265
$ws.Range("C6").Value = 10.26992913827406
$ws.Range("D6").Value = 6.417821099188971
$ws.Range("E6").Value = 12.61465101631297
$ws.Range("F6").Value = 33.0120193641583
$ws.Range("G6").Value = 44.08064239236445
$ws.Range("H6").Value = 17.98054735314724
$ws.Range("L6").Value = 9.143791646787019
$ws.Range("M6").Value = 53.63441320469072
$ws.Range("C7").Value = 10.21108469146095
$ws.Range("D7").Value = 6.369101455322975
$ws.Range("E7").Value = 12.68460138931299
$ws.Range("F7").Value = 32.86999723242284
$ws.Range("G7").Value = 43.85582930982749
$ws.Range("H7").Value = 17.8962745027927
$ws.Range("L7").Value = 9.200984226152123
$ws.Range("M7").Value = 54.65867870502068
$ws.Range("C8").Value = 9.964343339001132
$ws.Range("D8").Value = 6.161345250947249
$ws.Range("E8").Value = 13.00188707329523
$ws.Range("F8").Value = 32.31043888825692
$ws.Range("G8").Value = 42.98344954979327
$ws.Range("H8").Value = 17.54847232866545
$ws.Range("L8").Value = 9.455298144796219
$ws.Range("M8").Value = 58.95385739808107
$ws.Range("C9").Value = 9.528256118142124
$ws.Range("D9").Value = 5.779136271469241
$ws.Range("E9").Value = 13.65048903167136
$ws.Range("F9").Value = 31.47258834083744
$ws.Range("G9").Value = 41.74339042738088
$ws.Range("H9").Value = 16.95680567027191
$ws.Range("L9").Value = 9.957153005663578
$ws.Range("M9").Value = 66.54671658721102
$ws.Range("C10").Value = 9.2384935812825
$ws.Range("D10").Value = 5.512895966992443
$ws.Range("E10").Value = 14.1391946321805
$ws.Range("F10").Value = 31.0325902699491
$ws.Range("G10").Value = 41.15792422756581
$ws.Range("H10").Value = 16.58115358691739
$ws.Range("L10").Value = 10.3246284311248
$ws.Range("M10").Value = 71.60707953891347
$ws.Range("C11").Value = 9.113734966914166
$ws.Range("D11").Value = 5.394757676990878
$ws.Range("E11").Value = 14.36338855677676
$ws.Range("F11").Value = 30.87443711549043
$ws.Range("G11").Value = 40.97074616291753
$ws.Range("H11").Value = 16.42404279737925
$ws.Range("L11").Value = 10.49093293735867
$ws.Range("M11").Value = 73.79541396987965
$ws.Range("C12").Value = 9.067548095141877
$ws.Range("D12").Value = 5.350434306537792
$ws.Range("E12").Value = 14.44849541339263
$ws.Range("F12").Value = 30.82091682294665
$ws.Range("G12").Value = 40.91194346910247
$ws.Range("H12").Value = 16.36661824782135
$ws.Range("L12").Value = 10.55373904629919
$ws.Range("M12").Value = 74.60769295775361
$ws.Range("C13").Value = 9.077447544928189
$ws.Range("D13").Value = 5.359962011471662
$ws.Range("E13").Value = 14.43015772894482
$ws.Range("F13").Value = 30.83215485033414
$ws.Range("G13").Value = 40.92405958289878
$ws.Range("H13").Value = 16.37889210168809
$ws.Range("L13").Value = 10.54022086455803
$ws.Range("M13").Value = 74.43348376499792
$ws.Range("C14").Value = 9.109913710247794
$ws.Range("D14").Value = 5.391102976491348
$ws.Range("E14").Value = 14.37038648197145
$ws.Range("F14").Value = 30.86990443305824
$ws.Range("G14").Value = 40.96566243593838
$ws.Range("H14").Value = 16.41927642714008
$ws.Range("L14").Value = 10.49610368463369
$ws.Range("M14").Value = 73.86256940545707
$ws.Range("C15").Value = 9.129939106743587
$ws.Range("D15").Value = 5.410231053249563
$ws.Range("E15").Value = 14.33380042548198
$ws.Range("F15").Value = 30.8938663706561
$ws.Range("G15").Value = 40.9927386885584
$ws.Range("H15").Value = 16.44428531329789
$ws.Range("L15").Value = 10.46905723531486
$ws.Range("M15").Value = 73.51073178680443
$ws.Range("C16").Value = 9.246792278891421
$ws.Range("D16").Value = 5.520675205389236
$ws.Range("E16").Value = 14.12457554848846
$ws.Range("F16").Value = 31.04379919870747
$ws.Range("G16").Value = 41.17180919531785
$ws.Range("H16").Value = 16.59170579539568
$ws.Range("L16").Value = 10.31373868239429
$ws.Range("M16").Value = 71.46177432329418
$ws.Range("C17").Value = 9.320313061355739
$ws.Range("D17").Value = 5.589180901181531
$ws.Range("E17").Value = 13.99665951557341
$ws.Range("F17").Value = 31.14676905638782
$ws.Range("G17").Value = 41.30243382958592
$ws.Range("H17").Value = 16.68573286752104
$ws.Range("L17").Value = 10.21820109424296
$ws.Range("M17").Value = 70.17564559672547
$ws.Range("C18").Value = 9.363261553467156
$ws.Range("D18").Value = 5.628864587077145
$ws.Range("E18").Value = 13.9232658476092
$ws.Range("F18").Value = 31.20992825324403
$ws.Range("G18").Value = 41.38497464042958
$ws.Range("H18").Value = 16.74110560855052
$ws.Range("L18").Value = 10.16317176173844
$ws.Range("M18").Value = 69.42520135200165
$ws.Range("C19").Value = 9.3779155690656
$ws.Range("D19").Value = 5.642349509546719
$ws.Range("E19").Value = 13.89844889770033
$ws.Range("F19").Value = 31.23197852790991
$ws.Range("G19").Value = 41.41417269919467
$ws.Range("H19").Value = 16.76007297059039
$ws.Range("L19").Value = 10.14452773377745
$ws.Range("M19").Value = 69.16927890597209
$ws.Range("C20").Value = 9.312417923521449
$ws.Range("D20").Value = 5.581859373920532
$ws.Range("E20").Value = 14.0102581953262
$ws.Range("F20").Value = 31.1353982508931
$ws.Range("G20").Value = 41.28775674373316
$ws.Range("H20").Value = 16.67558930319491
$ws.Range("L20").Value = 10.22837971774357
$ws.Range("M20").Value = 70.31366381096564
$ws.Range("C21").Value = 9.100348573529363
$ws.Range("D21").Value = 5.381945032236829
$ws.Range("E21").Value = 14.38793752138181
$ws.Range("F21").Value = 30.8586409799218
$ws.Range("G21").Value = 40.9531094027022
$ws.Range("H21").Value = 16.40735765156258
$ws.Range("L21").Value = 10.50906692282947
$ws.Range("M21").Value = 74.03070599021051
$ws.Range("C22").Value = 8.967930969834697
$ws.Range("D22").Value = 5.253689561480914
$ws.Range("E22").Value = 14.6359691082283
$ws.Range("F22").Value = 30.71504107205583
$ws.Range("G22").Value = 40.80510765564448
$ws.Range("H22").Value = 16.24416398745831
$ws.Range("L22").Value = 10.69150423279482
$ws.Range("M22").Value = 76.36444365252592
$ws.Range("C23").Value = 9.038023660032017
$ws.Range("D23").Value = 5.321927398400505
$ws.Range("E23").Value = 14.50349933916902
$ws.Range("F23").Value = 30.78816287138499
$ws.Range("G23").Value = 40.87740277104246
$ws.Range("H23").Value = 16.33012411935729
$ws.Range("L23").Value = 10.59424030345171
$ws.Range("M23").Value = 75.12763537380624
$ws.Range("C24").Value = 9.315985199088782
$ws.Range("D24").Value = 5.585168504429992
$ws.Range("E24").Value = 14.00410976978474
$ws.Range("F24").Value = 31.14052668100163
$ws.Range("G24").Value = 41.29436912395298
$ws.Range("H24").Value = 16.68017111833773
$ws.Range("L24").Value = 10.22377828091217
$ws.Range("M24").Value = 70.25130013315504
$ws.Range("C25").Value = 9.640973731688449
$ws.Range("D25").Value = 5.87991615081001
$ws.Range("E25").Value = 13.47261139082614
$ws.Range("F25").Value = 31.66970891368799
$ws.Range("G25").Value = 42.02452624984006
$ws.Range("H25").Value = 17.10680378388997
$ws.Range("L25").Value = 9.821388876319199
$ws.Range("M25").Value = 64.5830220080871